$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7909169793128967
$ws.Range("B1").Value = 2.724210023880005
$ws.Range("C1").Value = 4.827226638793945
$ws.Range("D1").Value = 2.788937330245972
$ws.Range("E1").Value = 1.126551389694214
